# Insert 3 new data rows before row 1131 (shifting existing rows 1131:1197 down to 1134:1200)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1131:1133").Insert()

# Common values shared by the three new rows
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$fecha     = 44706
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100108
$producto  = "Tropicales y subtropicales"
$categoriaId = 100108002
$categoria = "Mango"
$variedad  = "Sin especificar"
$unidad    = "$/bandeja 4 kilos"
$origen    = "Brasil"
$kgUnidad  = 4

# Row 1131 - Calidad: Especial
$r = 1131
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 10000
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 10000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2500
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 1132 - Calidad: Primera
$r = 1132
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 10000
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 10000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2500
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 1133 - Calidad: Segunda
$r = 1133
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 10000
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 10000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2500
$ws.Cells.Item($r, 20).Value = $kgUnidad
